# Edit workbook per commit: "additional F_min expressions tested"
$wb = $excel.ActiveWorkbook

# --- Sheet: input_concentrations (sheet2) ---
$ws2 = $wb.Worksheets.Item("input_concentrations")

# Update A3 (new first data point) and the regenerated series in A4:B20
$ws2.Range("A3").Value = 0.00086623199999999998
$ws2.Range("A4").Value = 0.00086479200600440191
$ws2.Range("B4").Value = 0.000098894110120767572
$ws2.Range("A5").Value = 0.00086335679165531001
$ws2.Range("B5").Value = 0.00019745996964508995
$ws2.Range("A6").Value = 0.00086192633319515696
$ws2.Range("B6").Value = 0.00029569921016553907
$ws2.Range("A7").Value = 0.00086050060702356601
$ws2.Range("B7").Value = 0.00039361345247931315
$ws2.Range("A8").Value = 0.00085907958969605382
$ws2.Range("B8").Value = 0.00049120430667737303
$ws2.Range("A9").Value = 0.00085766325792274538
$ws2.Range("B9").Value = 0.00058847337223269752
$ws2.Range("A10").Value = 0.00085625158856710113
$ws2.Range("B10").Value = 0.0006854222380876661
$ws2.Range("A11").Value = 0.00085484455864465723
$ws2.Range("B11").Value = 0.00078205248274058057
$ws2.Range("A12").Value = 0.00085344214532177875
$ws2.Range("B12").Value = 0.00087836567433133565
$ws2.Range("A13").Value = 0.00085204432591442364
$ws2.Range("B13").Value = 0.00097436337072624702
$ws2.Range("A14").Value = 0.00085065107788692012
$ws2.Range("B14").Value = 0.0010700471196020466
$ws2.Range("A15").Value = 0.00084926237885075582
$ws2.Range("B15").Value = 0.0011654184585290557
$ws2.Range("A16").Value = 0.00084787820656337767
$ws2.Range("B16").Value = 0.0012604789150535466
$ws2.Range("A17").Value = 0.00084649853892700545
$ws2.Range("B17").Value = 0.0013552300067792952
$ws2.Range("A18").Value = 0.00084512335398745527
$ws2.Range("B18").Value = 0.0014496732414483447
$ws2.Range("A19").Value = 0.00084375262993297495
$ws2.Range("B19").Value = 0.001543810117020978
$ws2.Range("A20").Value = 0.00084238634509309113
$ws2.Range("B20").Value = 0.0016376421217549177

# Remove the now-unused trailing rows (21 and 22) - series shortened
$ws2.Rows("21:22").Delete()

# New bestFit-style helper column metadata (col G) introduced by the rerun
$ws2.Columns(7).ColumnWidth = 11.166666666666666

# --- Sheet: heats (sheet5) ---
$ws5 = $wb.Worksheets.Item("heats")

# Recomputed volumes (row 2) and observations (row 3) for columns B:R
$ws5.Range("B2").Value = 15.024977
$ws5.Range("C2").Value = 15.049954
$ws5.Range("D2").Value = 15.074930999999999
$ws5.Range("E2").Value = 15.099907999999999
$ws5.Range("F2").Value = 15.124885000000001
$ws5.Range("G2").Value = 15.149862000000001
$ws5.Range("H2").Value = 15.174839
$ws5.Range("I2").Value = 15.199816
$ws5.Range("J2").Value = 15.224793
$ws5.Range("K2").Value = 15.24977
$ws5.Range("L2").Value = 15.274747
$ws5.Range("M2").Value = 15.299723999999999
$ws5.Range("N2").Value = 15.324700999999999
$ws5.Range("O2").Value = 15.349678000000001
$ws5.Range("P2").Value = 15.374655000000001
$ws5.Range("Q2").Value = 15.399632
$ws5.Range("R2").Value = 15.424609
$ws5.Range("B3").Value = 0.073663000000000006
$ws5.Range("C3").Value = 0.074034000000000003
$ws5.Range("D3").Value = 0.074881000000000003
$ws5.Range("E3").Value = 0.073269000000000001
$ws5.Range("F3").Value = 0.074992000000000003
$ws5.Range("G3").Value = 0.074681999999999998
$ws5.Range("H3").Value = 0.071710999999999997
$ws5.Range("I3").Value = 0.064866999999999994
$ws5.Range("J3").Value = 0.041674000000000003
$ws5.Range("K3").Value = 0.019753
$ws5.Range("L3").Value = 0.0068329999999999997
$ws5.Range("M3").Value = 0.0014450000000000001
$ws5.Range("N3").Value = 0.0019889999999999999
$ws5.Range("O3").Value = 0.001794
$ws5.Range("P3").Value = -0.0068100000000000001
$ws5.Range("Q3").Value = 0.00030299999999999999
$ws5.Range("R3").Value = 0.00011

# Drop the last two simulated points (columns S and T)
$ws5.Columns("S:T").Delete()

# --- Sheet: component_name (sheet4) ---
$ws4 = $wb.Worksheets.Item("component_name")
$ws4.Range("A1").Value = "T3H"

# --- Selection / active-sheet bookkeeping ---
# Restore prior selections on sheets that are no longer active
$ws2.Range("L19").Select()
$ws5.Range("K12").Select()

# component_name becomes the active sheet/tab, selection on B1
$ws4.Range("B1").Select()
